# Kompose presentation — "adding logo for abgabe"
#
# 1. Slide 1: punch up the subtitle into a single bold title-style line
#    ("A Distributed Playlist for Android" -> "THE DISTRIBUTED PLAYLIST"),
#    bump the font size and turn off auto-shrink-to-fit so it sits at its
#    authored size like a logo lockup.
# 2. Slide 4: drop the redundant "And much more !!" title placeholder that
#    duplicated the body copy below it.

$p = $ppt.ActivePresentation

# --- Slide 1: subtitle / hero line -----------------------------------
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Text = "THE DISTRIBUTED PLAYLIST"
$subtitle.TextFrame.TextRange.Font.Size = 43
$subtitle.TextFrame.AutoSize = 0

# --- Slide 4: remove the now-redundant title placeholder --------------
$s4 = $p.Slides.Item(4)
$title4 = $s4.Shapes.Item("Titel 1")
$title4.Cut()
